# Sort the comma-separated "Recorded By" names in column G alphabetically
# using an ordinal (case-sensitive, byte-value) comparison - matching the
# sort order used by the upstream sync job.

function Compare-Ordinal($s1, $s2) {
    $cmpLen1 = $s1.Length
    $cmpLen2 = $s2.Length
    $cmpMinLen = [Math]::Min($cmpLen1, $cmpLen2)
    for ($cmpIdx = 0; $cmpIdx -lt $cmpMinLen; $cmpIdx++) {
        $cc1 = [int][char]$s1[$cmpIdx]
        $cc2 = [int][char]$s2[$cmpIdx]
        if ($cc1 -lt $cc2) { return -1 }
        if ($cc1 -gt $cc2) { return 1 }
    }
    if ($cmpLen1 -lt $cmpLen2) { return -1 }
    if ($cmpLen1 -gt $cmpLen2) { return 1 }
    return 0
}

function Sort-Ordinal($items) {
    $arr = @($items)
    $n = $arr.Length
    for ($si = 1; $si -lt $n; $si++) {
        $key = $arr[$si]
        $sj = $si - 1
        while ($sj -ge 0 -and (Compare-Ordinal $arr[$sj] $key) -gt 0) {
            $arr[$sj + 1] = $arr[$sj]
            $sj = $sj - 1
        }
        $arr[$sj + 1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $rawParts = $val.Split(",")
        $parts = @()
        foreach ($p in $rawParts) {
            $parts += $p.Trim()
        }
        if ($parts.Length -gt 1) {
            $sortedParts = Sort-Ordinal $parts
            $joined = [string]::Join(", ", $sortedParts)
            if ($joined -ne $val) {
                $cell.Value2 = $joined
            }
        }
    }
}
